$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet currently ends at row 148, which holds the "Maximum Average Pass
# Ratio" entry (dated 2025-09-01 / serial 45901). We need to:
#   1) push that entry down to row 151
#   2) add two "date stamp only" rows (148, 149) for 2025-08-30 and 2025-08-31
#   3) add a brand new row (150) for LeetCode 37 "Sudoku Solver"
# ---------------------------------------------------------------------------

# Remember the values currently sitting in row 148 before we move anything.
$oldA = $ws.Range("A148").Value()
$oldB = $ws.Range("B148").Value()
$oldC = $ws.Range("C148").Value()
$oldD = $ws.Range("D148").Value()
$oldE = $ws.Range("E148").Value()
$oldF = $ws.Range("F148").Value()
$oldG = $ws.Range("G148").Value()
$oldH = $ws.Range("H148").Value()
$oldI = $ws.Range("I148").Value()

# Drop row 148 entirely so the sheet is back to 147 rows with no leftover
# row-level formatting (height, etc.) to clean up afterwards.
$ws.Rows.Item(148).Delete()

# --- Row 148: just the daily H/I date stamps (2025-08-30) -----------------
$ws.Range("H147:I147").Copy()
$ws.Range("H148:I148").PasteSpecial(-4122)
$ws.Range("H148").Value = 45899
$ws.Range("I148").Value = 45899

# --- Row 149: just the daily H/I date stamps (2025-08-31) -----------------
$ws.Range("H147:I147").Copy()
$ws.Range("H149:I149").PasteSpecial(-4122)
$ws.Range("H149").Value = 45900
$ws.Range("I149").Value = 45900

# --- Row 150: new entry, LeetCode 37 "Sudoku Solver" -----------------------
$ws.Range("A147:F147").Copy()
$ws.Range("A150:F150").PasteSpecial(-4122)
$ws.Range("H147:I147").Copy()
$ws.Range("H150:I150").PasteSpecial(-4122)

$ws.Range("A150").Value = 37
$ws.Range("B150").Value = "Sudoku Solver"
$ws.Range("C150").Value = "#matrix #bit-minipulation #hash-table #backtracking "
$ws.Range("D150").Value = "hard"
$ws.Range("E150").Value = 0
$ws.Range("F150").Value = 1
$ws.Range("H150").Value = 45900
$ws.Range("I150").Value = 45900
$ws.Rows.Item(150).RowHeight = 68

# --- Row 151: the entry that used to be row 148 ----------------------------
$ws.Range("A147:I147").Copy()
$ws.Range("A151:I151").PasteSpecial(-4122)

$ws.Range("A151").Value = $oldA
$ws.Range("B151").Value = $oldB
$ws.Range("C151").Value = $oldC
$ws.Range("D151").Value = $oldD
$ws.Range("E151").Value = $oldE
$ws.Range("F151").Value = $oldF
$ws.Range("G151").Value = $oldG
$ws.Range("H151").Value = $oldH
$ws.Range("I151").Value = $oldI
$ws.Rows.Item(151).RowHeight = 34

# Match the author's final cursor position.
$null = $ws.Range("H155").Select()
